$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-117 down to 70-118.
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new Frambuesa price record.
$ws.Cells.Item(69,1).Value2  = 9
$ws.Cells.Item(69,2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(69,3).Value2  = "Metropolitana"
$ws.Cells.Item(69,4).Value2  = 44942
$ws.Cells.Item(69,5).Value2  = 13
$ws.Cells.Item(69,6).Value2  = "Fruta"
$ws.Cells.Item(69,7).Value2  = 100101
$ws.Cells.Item(69,8).Value2  = "Berries"
$ws.Cells.Item(69,9).Value2  = 100101004
$ws.Cells.Item(69,10).Value2 = "Frambuesa"
$ws.Cells.Item(69,11).Value2 = "Sin especificar"
$ws.Cells.Item(69,12).Value2 = "Primera"
$ws.Cells.Item(69,13).Value2 = 300
$ws.Cells.Item(69,14).Value2 = 8000
$ws.Cells.Item(69,15).Value2 = 8000
$ws.Cells.Item(69,16).Value2 = 8000
$ws.Cells.Item(69,17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(69,18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(69,19).Value2 = 4000
$ws.Cells.Item(69,20).Value2 = 2
